$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: set B-column numeric/text values and A-column cells that reuse existing shared strings ---
# (write order here does not affect final shared-string table position for already-present strings)
$ws.Range("A1").Value = "Labels"
$ws.Range("B1").Value = "Values"
$ws.Range("A2").Value = "Congress"
$ws.Range("B2").Value = 105
$ws.Range("A3").Value = "Session"
$ws.Range("B3").Value = 2
$ws.Range("A4").Value = "Start Date"
$ws.Range("B4").Value = 35822
$ws.Range("A5").Value = "End Date"
$ws.Range("B5").Value = 36089
$ws.Range("A6").Value = "Civilian "
$ws.Range("A17").Value = "Air Force "
$ws.Range("A22").Value = "Army "
$ws.Range("A27").Value = "Navy "
$ws.Range("A32").Value = "Marine Corps"

# --- Step 2: set B-column values for rows introducing brand-new shared strings (values first, harmless) ---
$ws.Range("B7").Value = 336
$ws.Range("B8").Value = 124
$ws.Range("B9").Value = 319
$ws.Range("B10").Value = 24
$ws.Range("B11").Value = 114
$ws.Range("B18").Value = 6070
$ws.Range("B19").Value = 21
$ws.Range("B20").Value = 6087
$ws.Range("B21").Value = 4
$ws.Range("B23").Value = 5479
$ws.Range("B24").Value = 2
$ws.Range("B25").Value = 5478
$ws.Range("B26").Value = 3
$ws.Range("B28").Value = 5047
$ws.Range("B29").Value = 4
$ws.Range("B30").Value = 5045
$ws.Range("B31").Value = 6
$ws.Range("B33").Value = 1847
$ws.Range("B34").Value = 1847
$ws.Range("B35").Value = 20225
$ws.Range("B36").Value = 237
$ws.Range("B37").Value = 20302
$ws.Range("B38").Value = 27
$ws.Range("B39").Value = 133
$ws.Range("B13").Value = 1446
$ws.Range("B14").Value = 86
$ws.Range("B15").Value = 1526
$ws.Range("B16").Value = 6

# --- Step 3: set A-column cells that introduce brand-new shared strings, in the exact order they
#     must first appear so the rebuilt shared-strings table matches the target order ---
$ws.Range("A7").Value = "     Civilian, New nominations"
$ws.Range("A8").Value = "     Civilian, Carryover nominations"
$ws.Range("A9").Value = "     Civilian, Confirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("A11").Value = "     Civilian, Returned to White House "
$ws.Range("A18").Value = "     Air Force, New nominations"
$ws.Range("A19").Value = "     Air Force, Carryover nominations"
$ws.Range("A20").Value = "     Air Force, Confirmed "
$ws.Range("A21").Value = "     Air Force, Returned to White House "
$ws.Range("A23").Value = "     Army, New nominations"
$ws.Range("A24").Value = "     Army, Carryover nominations"
$ws.Range("A25").Value = "     Army, Confirmed "
$ws.Range("A26").Value = "     Army, Returned to White House "
$ws.Range("A28").Value = "     Navy, New nominations"
$ws.Range("A29").Value = "     Navy, Carryover nominations"
$ws.Range("A30").Value = "     Navy, Confirmed "
$ws.Range("A31").Value = "     Navy, Returned to White House "
$ws.Range("A33").Value = "     Marine Corps, New nominations"
$ws.Range("A34").Value = "     Marine Corps, Confirmed "
$ws.Range("A35").Value = "Total new nominations"
$ws.Range("A36").Value = "Total carryover nominations"
$ws.Range("A37").Value = "Total confirmed "
$ws.Range("A38").Value = "Total withdrawn "
$ws.Range("A39").Value = "Total returned"
$ws.Range("A12").Value = "Civilian (FS, PHS, CG, NOAA)"
$ws.Range("A13").Value = "     Civilian (FS, PHS, CG, NOAA), New nominations"
$ws.Range("A14").Value = "     Civilian (FS, PHS, CG, NOAA), Carryover nominations"
$ws.Range("A15").Value = "     Civilian (FS, PHS, CG, NOAA), Confirmed "
$ws.Range("A16").Value = "     Civilian (FS, PHS, CG, NOAA), Returned to White House "

# --- Step 4: new B35 cell needs the same "#,##0" number format / right alignment as the other Total rows ---
$ws.Range("B35").NumberFormat = $ws.Range("B34").NumberFormat
$ws.Range("B35").HorizontalAlignment = $ws.Range("B34").HorizontalAlignment

# --- Step 5: drop the old trailing row 40 (sheet now spans only to row 39) ---
$ws.Rows(40).Delete()
